# Append a new row (row 73) of database records to each of the four
# worksheets, mirroring the existing row layout (time / hex-string columns /
# numeric decoded columns).

$wb = $excel.ActiveWorkbook

$rowsData = @{
    "DE_LFT_#1" = @{
        A = 45859.43543981481
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x44"
        E = "0x14"
        F = 380
        G = ("7.598631275147109e+23" -as [double])
        H = 324
        I = 14
    }
    "DE_LFT_#2" = @{
        A = 45859.43543981481
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x48"
        E = "0xe"
        F = 380
        G = ("5.68432987514711e+23" -as [double])
        H = 328
        I = 14
    }
    "DE_PLT_#1" = @{
        A = 45859.43543981481
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x79"
        E = "0x7"
        F = 130
        G = ("5.68631262647114e+23" -as [double])
        H = 121
        I = 7
    }
    "DE_PLT_#2" = @{
        A = 45859.43543981481
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x79"
        E = "0x3"
        F = 130
        G = ("9.85046333984776e+23" -as [double])
        H = 121
        I = 3
    }
}

foreach ($sheetName in $rowsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 73

    $data = $rowsData[$sheetName]

    # Column A: numeric timestamp, same style/number format as the row above.
    $cellA = $ws.Cells.Item($newRow, 1)
    $cellA.Value = $data.A
    $cellA.NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    # Columns B-E: text hex strings.
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    # Columns F-I: decoded numeric values.
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
